# Add the two new opponent rows (Alabama, Michigan State) to the bottom of the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Alabama"
$ws.Range("B4").Value = "https://www.sports-reference.com/cbb/schools/alabama/2023.html"
$ws.Range("A5").Value = "Michigan State"
$ws.Range("B5").Value = "https://www.sports-reference.com/cbb/schools/michigan-state/2023.html"

# Column A needs to widen to fit the longest new entry ("Michigan State"),
# mirroring Excel's "best fit" autosize that happens when the data changes.
$ws.Columns("A:A").AutoFit() | Out-Null
$ws.Columns("A:A").ColumnWidth = 12.25

# Leave the selection on the last cell that was touched, matching the saved view.
$ws.Range("B5").Select() | Out-Null
